# BUG: Removed dumping of sex and species variable
# Adds new checklist columns K:O to row 1 (headers), and the corresponding
# "1" markers to rows 4, 5, 6, 7, 8 that were previously missing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in row 1
$ws.Range("K1").Value = "help"
$ws.Range("L1").Value = "help"
$ws.Range("M1").Value = "test.pkl"
$ws.Range("N1").Value = "save"
$ws.Range("O1").Value = "sa"

# New data markers
$ws.Range("K4").Value = 1
$ws.Range("O5").Value = 1
$ws.Range("L6").Value = 1
$ws.Range("N7").Value = 1
$ws.Range("M8").Value = 1
